$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2558
$ws.Range("C2").Value = 2514.2
$ws.Range("D2").Value = 2537.9
$ws.Range("E2").Value = 2540.05
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 2521.6

$ws.Range("B3").Value = 391.1
$ws.Range("C3").Value = 375.3
$ws.Range("D3").Value = 376.5
$ws.Range("E3").Value = 376.4
$ws.Range("F3").Value = 45
$ws.Range("G3").Value = 386.6

$ws.Range("B4").Value = 1567.6
$ws.Range("C4").Value = 1531.3
$ws.Range("D4").Value = 1538.45
$ws.Range("E4").Value = 1538.1
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = 1549.75

$ws.Range("B5").Value = 7565.15
$ws.Range("C5").Value = 7455.55
$ws.Range("D5").Value = 7491.5
$ws.Range("E5").Value = 7492.2
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 7483.45

$ws.Range("B6").Value = 250.2
$ws.Range("C6").Value = 240.5
$ws.Range("D6").Value = 249.25
$ws.Range("E6").Value = 249.35
$ws.Range("F6").Value = 203
$ws.Range("G6").Value = 241.25

$ws.Range("B7").Value = 212.65
$ws.Range("C7").Value = 209.25
$ws.Range("D7").Value = 210.5
$ws.Range("E7").Value = 210.8
$ws.Range("F7").Value = 162
$ws.Range("G7").Value = 210.95

$ws.Range("B8").Value = 284.15
$ws.Range("C8").Value = 277.55
$ws.Range("D8").Value = 278.7
$ws.Range("E8").Value = 279.3
$ws.Range("F8").Value = 132
$ws.Range("G8").Value = 282.75

$ws.Range("B9").Value = 536.55
$ws.Range("C9").Value = 528.05
$ws.Range("D9").Value = 530.25
$ws.Range("E9").Value = 532.4
$ws.Range("F9").Value = 23
$ws.Range("G9").Value = 533.9

$ws.Range("B10").Value = 3447.9
$ws.Range("C10").Value = 3390.9
$ws.Range("D10").Value = 3422
$ws.Range("E10").Value = 3427.2
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 3399.1

$ws.Range("B11").Value = 148.8
$ws.Range("C11").Value = 146.65
$ws.Range("D11").Value = 148.2
$ws.Range("E11").Value = 148.25
$ws.Range("F11").Value = 86
$ws.Range("G11").Value = 147.55

$ws.Range("B12").Value = 1310
$ws.Range("C12").Value = 1296.2
$ws.Range("D12").Value = 1304.15
$ws.Range("E12").Value = 1304.6
$ws.Range("F12").Value = 42
$ws.Range("G12").Value = 1301.8

$ws.Range("B13").Value = 1670
$ws.Range("C13").Value = 1652.6
$ws.Range("D13").Value = 1658.9
$ws.Range("E13").Value = 1661.75
$ws.Range("F13").Value = 867
$ws.Range("G13").Value = 1657.5

$ws.Range("B14").Value = 502.85
$ws.Range("C14").Value = 494.95
$ws.Range("D14").Value = 496.25
$ws.Range("E14").Value = 496.55
$ws.Range("F14").Value = 85
$ws.Range("G14").Value = 502.2

$ws.Range("B15").Value = 997.5
$ws.Range("C15").Value = 988.6
$ws.Range("D15").Value = 991.6
$ws.Range("E15").Value = 992.45
$ws.Range("F15").Value = 133
$ws.Range("G15").Value = 991.85

$ws.Range("B16").Value = 1456
$ws.Range("C16").Value = 1446.45
$ws.Range("D16").Value = 1451.95
$ws.Range("E16").Value = 1450
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 1453

$ws.Range("B17").Value = 1517.9
$ws.Range("C17").Value = 1509
$ws.Range("D17").Value = 1512
$ws.Range("E17").Value = 1511.6
$ws.Range("F17").Value = 86
$ws.Range("G17").Value = 1515

$ws.Range("B18").Value = 718.2
$ws.Range("C18").Value = 708.75
$ws.Range("D18").Value = 710.35
$ws.Range("E18").Value = 711.45
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 715.55

$ws.Range("B19").Value = 465.9
$ws.Range("C19").Value = 457.45
$ws.Range("D19").Value = 459
$ws.Range("E19").Value = 458.75
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 461.75

$ws.Range("B20").Value = 1610.7
$ws.Range("C20").Value = 1571.85
$ws.Range("D20").Value = 1600.05
$ws.Range("E20").Value = 1601.1
$ws.Range("F20").Value = 54
$ws.Range("G20").Value = 1586.4

$ws.Range("B21").Value = 305.45
$ws.Range("C21").Value = 299.95
$ws.Range("D21").Value = 301.05
$ws.Range("E21").Value = 300.95
$ws.Range("F21").Value = 25
$ws.Range("G21").Value = 302.6

$ws.Range("B22").Value = 2467
$ws.Range("C22").Value = 2449.5
$ws.Range("D22").Value = 2452.5
$ws.Range("E22").Value = 2457.85
$ws.Range("F22").Value = 127
$ws.Range("G22").Value = 2462.5

$ws.Range("B23").Value = 600.7
$ws.Range("C23").Value = 595.2
$ws.Range("D23").Value = 598.05
$ws.Range("E23").Value = 598.8
$ws.Range("F23").Value = 117
$ws.Range("G23").Value = 597.6

$ws.Range("B24").Value = 601.8
$ws.Range("C24").Value = 589.55
$ws.Range("D24").Value = 597
$ws.Range("E24").Value = 596.2
$ws.Range("F24").Value = 19
$ws.Range("G24").Value = 598.35

$ws.Range("B25").Value = 1074.1
$ws.Range("C25").Value = 1063.35
$ws.Range("D25").Value = 1068.7
$ws.Range("E25").Value = 1070.45
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 1071.1

$ws.Range("B26").Value = 636.5
$ws.Range("C26").Value = 627.15
$ws.Range("D26").Value = 634.2
$ws.Range("E26").Value = 634.25
$ws.Range("F26").Value = 189
$ws.Range("G26").Value = 631.75

$ws.Range("B27").Value = 267.4
$ws.Range("C27").Value = 263.1
$ws.Range("D27").Value = 264.25
$ws.Range("E27").Value = 263.95
$ws.Range("F27").Value = 76
$ws.Range("G27").Value = 266.85

$ws.Range("B28").Value = 132.5
$ws.Range("C28").Value = 129.85
$ws.Range("D28").Value = 132.2
$ws.Range("E28").Value = 131.95
$ws.Range("F28").Value = 850
$ws.Range("G28").Value = 132.3

$ws.Range("B29").Value = 8750.950000000001
$ws.Range("C29").Value = 8661.75
$ws.Range("D29").Value = 8731.049999999999
$ws.Range("E29").Value = 8728.35
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 8682.75
